# Apply the "5-year growth is frozen before first drop in projections" update
# to the GroupGrowthRateMax sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GroupGrowthRateMax")

# --- 1. Update the "rate" (column D) values: for each 6-row group (one per
#     region/group_name, years 2025/2030/2035/2040/2045/2050), growth is
#     frozen at the value immediately preceding the first year-over-year drop.
$rateRanges = @(
    @{ Range = "D3:D7";   Value = 12.93494652 },
    @{ Range = "D9:D13";  Value = 13.391031024 },
    @{ Range = "D15:D19"; Value = 13.435053768 },
    @{ Range = "D22:D25"; Value = 14.60408571 },
    @{ Range = "D28:D31"; Value = 25.636466032 },
    @{ Range = "D34:D37"; Value = 25.636466032 },
    @{ Range = "D39:D43"; Value = 13.8333099 },
    @{ Range = "D46:D49"; Value = 25.212515166 },
    @{ Range = "D52:D55"; Value = 26.012634867 },
    @{ Range = "D58:D61"; Value = 26.012634867 },
    @{ Range = "D64:D67"; Value = 26.012634867 },
    @{ Range = "D70:D73"; Value = 24.811278945 },
    @{ Range = "D76:D79"; Value = 24.811278945 },
    @{ Range = "D82:D85"; Value = 24.811278945 },
    @{ Range = "D88:D91"; Value = 24.811278945 },
    @{ Range = "D94:D97"; Value = 25.367137145 },
    @{ Range = "D100:D103"; Value = 25.367137145 },
    @{ Range = "D106:D109"; Value = 25.367137145 },
    @{ Range = "D112:D115"; Value = 25.367137145 },
    @{ Range = "D117:D121"; Value = 13.289390928 },
    @{ Range = "D124:D127"; Value = 21.879636292 },
    @{ Range = "D130:D133"; Value = 25.750595548 },
    @{ Range = "D136:D139"; Value = 25.992823061 },
    @{ Range = "D142:D145"; Value = 25.992823061 }
)

foreach ($item in $rateRanges) {
    $ws.Range($item.Range).Value = $item.Value
}

# --- 2. Update the "notes" (column E) text for every data row (2-145):
#     insert ", 5-year growth is frozen before first drop in projections"
#     right before the "; 90th percentile" clause.
$oldFragment = "; 90th percentile"
$newFragment = ", 5-year growth is frozen before first drop in projections; 90th percentile"

$lastRow = 145
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("E$r")
    $text = $cell.Value2
    if ($text -ne $null -and $text.ToString().Contains($oldFragment)) {
        $cell.Value2 = $text.ToString().Replace($oldFragment, $newFragment)
    }
}
